# Added more simulated games, which changed the relative transition-probability
# matrix produced by the (sped-up) simulate-game logic. Update the affected
# probabilities in-place on the "UC Merced_B" transition matrix sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2258064516129032
$ws.Range("C2").Value = 0.4838709677419355
$ws.Range("J2").Value = 0.03225806451612903
$ws.Range("O2").Value = 0.03225806451612903
$ws.Range("P2").Value = 0.1612903225806452
$ws.Range("S2").Value = 0.06451612903225806
$ws.Range("J3").Value = 0.1333333333333333
$ws.Range("P3").Value = 0.6666666666666666
$ws.Range("S3").Value = 0.2
$ws.Range("J4").Value = 0.25
$ws.Range("P4").Value = 0.5
$ws.Range("S4").Value = 0.25
$ws.Range("D6").Value = 0.06666666666666667
$ws.Range("J6").Value = 0.4
$ws.Range("Q6").Value = 0.1333333333333333
$ws.Range("R6").Value = 0.1333333333333333
$ws.Range("S6").Value = 0.2666666666666667
$ws.Range("J7").Value = 0.08333333333333333
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("B8").Value = 0.1111111111111111
$ws.Range("D8").Value = 0.03703703703703703
$ws.Range("F8").Value = 0.03703703703703703
$ws.Range("J8").Value = 0.2222222222222222
$ws.Range("O8").Value = 0.03703703703703703
$ws.Range("Q8").Value = 0.07407407407407407
$ws.Range("R8").Value = 0.1111111111111111
$ws.Range("S8").Value = 0.3703703703703703
$ws.Range("J9").Value = 0.3636363636363636
$ws.Range("Q9").Value = 0.09090909090909091
$ws.Range("R9").Value = 0.2727272727272727
$ws.Range("S9").Value = 0.2727272727272727
$ws.Range("B10").Value = 0.1651376146788991
$ws.Range("D10").Value = 0.01834862385321101
$ws.Range("F10").Value = 0.06422018348623854
$ws.Range("J10").Value = 0.1651376146788991
$ws.Range("O10").Value = 0.03669724770642202
$ws.Range("Q10").Value = 0.1192660550458716
$ws.Range("R10").Value = 0.1192660550458716
$ws.Range("S10").Value = 0.3119266055045872
$ws.Range("G11").Value = 0.0625
$ws.Range("J11").Value = 0.0625
$ws.Range("K11").Value = 0.1875
$ws.Range("L11").Value = 0.6875
$ws.Range("G12").Value = 0.4615384615384616
$ws.Range("J12").Value = 0.3846153846153846
$ws.Range("L12").Value = 0.1538461538461539
$ws.Range("G14").Value = 1
$ws.Range("H15").Value = 0.1333333333333333
$ws.Range("I15").Value = 0.1333333333333333
$ws.Range("J15").Value = 0.4
$ws.Range("S15").Value = 0.3333333333333333
$ws.Range("H16").Value = 0.2222222222222222
$ws.Range("J16").Value = 0.5555555555555556
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.05555555555555555
$ws.Range("F17").Value = 0.1
$ws.Range("H17").Value = 0.2
$ws.Range("I17").Value = 0.1
$ws.Range("J17").Value = 0.4
$ws.Range("K17").Value = 0.1
$ws.Range("S17").Value = 0.1
$ws.Range("H18").Value = 0.2608695652173913
$ws.Range("J18").Value = 0.5217391304347826
$ws.Range("K18").Value = 0.04347826086956522
$ws.Range("S18").Value = 0.1739130434782609
$ws.Range("F19").Value = 0.03703703703703703
$ws.Range("H19").Value = 0.1358024691358025
$ws.Range("I19").Value = 0.08641975308641975
$ws.Range("J19").Value = 0.3580246913580247
$ws.Range("K19").Value = 0.09876543209876543
$ws.Range("M19").Value = 0.04938271604938271
$ws.Range("N19").Value = 0.01234567901234568
$ws.Range("O19").Value = 0.08641975308641975
$ws.Range("S19").Value = 0.1358024691358025
